$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update effort/remaining values for row 5 (Discord Library learning)
$ws.Range("D5").Value = "2h"
$ws.Range("E5").Value = "-"

# Update the active cell selection
$ws.Range("E6").Select()
